# Add beneficiary fields (telephone, gender, address, etc.) to the bulk
# upload template, on both the blank "Sheet1" template and the populated
# "Example" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Example")

# --- Sheet1: header row only -------------------------------------------
$ws1.Range("O1").Value = "telephone"
$ws1.Range("P1").Value = "gender"
$ws1.Range("Q1").Value = "addressLine1"
$ws1.Range("R1").Value = "addressLine2"
$ws1.Range("S1").Value = "city"
$ws1.Range("T1").Value = "stateOrProvince"
$ws1.Range("U1").Value = "postalCode"
$ws1.Range("V1").Value = "Country"
$null = $ws1.Range("O1:V1").Select()

# --- Example sheet: header row (same field names, strings reused) ------
$ws2.Range("O1").Value = "telephone"
$ws2.Range("P1").Value = "gender"
$ws2.Range("Q1").Value = "addressLine1"
$ws2.Range("R1").Value = "addressLine2"
$ws2.Range("S1").Value = "city"
$ws2.Range("T1").Value = "stateOrProvince"
$ws2.Range("U1").Value = "postalCode"
$ws2.Range("V1").Value = "Country"

# --- Example sheet: sample data row -------------------------------------
$ws2.Range("O2").Value = "91 (80) 26677444"
$ws2.Range("P2").Value = "FEMALE"

# --- Example sheet: gender "Choices:" list, to the right of the other
# choice lists (fundRaisingReason/fundRaisingFor/category) --------------
$ws2.Range("P5").Value = "Choices:"
$ws2.Range("P5").Font.Bold = $true
$ws2.Range("P6").Value = "MALE"
$ws2.Range("P7").Value = "FEMALE"
$ws2.Range("P8").Value = "UNSPECIFIED"

# --- Example sheet: remaining sample data for the new address fields ---
$ws2.Range("Q2").Value = "123 Sample Street"
$ws2.Range("R2").Value = "Basavanagudi"
$ws2.Range("S2").Value = "Bangalore"
$ws2.Range("T2").Value = "Karnataka"
$ws2.Range("U2").Value = 560028
$ws2.Range("V2").Value = "India"

$null = $ws2.Range("A1:V2").Select()

# Column O on the Example sheet needs to be a bit wider to fit the phone
# number sample.
$ws2.Columns.Item(15).ColumnWidth = 12.1640625

# Give the Example sheet an explicit portrait page setup (Sheet1 already
# has one).
$ws2.PageSetup.Orientation = 1
